# Update countries & provincias Spain
# - Refresh the "datos actualizados" timestamp string.
# - Update several countries' case counts (Estados Unidos, Reino Unido, Alemania,
#   India, Peru, Finlandia, Oman).
# - Estado de Palestina's case count dropped, pushing it below Congo (and the block
#   Sierra Leona..Congo each shift up a row); re-write rows 127-135 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp string (A1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 20 de Mayo de 2020 a las 22:05"

# Estados Unidos (row 4)
$ws.Cells.Item(4, 2).Value = 1582373
$ws.Cells.Item(4, 3).Value = 11790
$ws.Cells.Item(4, 5).Value = 1122299
$ws.Cells.Item(4, 7).Value = 794
$ws.Cells.Item(4, 8).Value = 94327

# Reino Unido (row 8)
$ws.Cells.Item(8, 2).Value = 248293

# Alemania (row 11)
$ws.Cells.Item(11, 2).Value = 178443
$ws.Cells.Item(11, 3).Value = 616
$ws.Cells.Item(11, 5).Value = 13284
$ws.Cells.Item(11, 7).Value = 66
$ws.Cells.Item(11, 8).Value = 8259

# India (row 14)
$ws.Cells.Item(14, 2).Value = 112028
$ws.Cells.Item(14, 3).Value = 5553
$ws.Cells.Item(14, 5).Value = 63172

# Peru (row 15)
$ws.Cells.Item(15, 2).Value = 104020
$ws.Cells.Item(15, 3).Value = 4537
$ws.Cells.Item(15, 4).Value = 41968
$ws.Cells.Item(15, 5).Value = 59028
$ws.Cells.Item(15, 7).Value = 110
$ws.Cells.Item(15, 8).Value = 3024

# Finlandia (row 62)
$ws.Cells.Item(62, 4).Value = 4800
$ws.Cells.Item(62, 5).Value = 1339

# Oman (row 65)
$ws.Cells.Item(65, 5).Value = 4352
$ws.Cells.Item(65, 7).Value = 3
$ws.Cells.Item(65, 8).Value = 30

# Rows 127-135: Sierra Leona..Congo each move up one row (unchanged totals);
# Estado de Palestina's total cases fell (577 -> 398) so it now lands on row 135,
# just above Etiopia (row 136, untouched).
$ws.Cells.Item(127, 1).Value = "Sierra Leona"
$ws.Cells.Item(127, 2).Value = 570
$ws.Cells.Item(127, 3).Value = 36
$ws.Cells.Item(127, 4).Value = 205
$ws.Cells.Item(127, 5).Value = 331
$ws.Cells.Item(127, 6).Value = 0
$ws.Cells.Item(127, 7).Value = 1
$ws.Cells.Item(127, 8).Value = 34

$ws.Cells.Item(128, 1).Value = "Republica del Chad"
$ws.Cells.Item(128, 2).Value = 565
$ws.Cells.Item(128, 3).Value = 20
$ws.Cells.Item(128, 4).Value = 177
$ws.Cells.Item(128, 5).Value = 331
$ws.Cells.Item(128, 6).Value = 0
$ws.Cells.Item(128, 7).Value = 1
$ws.Cells.Item(128, 8).Value = 57

$ws.Cells.Item(129, 1).Value = "Jamaica"
$ws.Cells.Item(129, 2).Value = 520
$ws.Cells.Item(129, 3).Value = 0
$ws.Cells.Item(129, 4).Value = 145
$ws.Cells.Item(129, 5).Value = 366
$ws.Cells.Item(129, 6).Value = 0
$ws.Cells.Item(129, 7).Value = 0
$ws.Cells.Item(129, 8).Value = 9

$ws.Cells.Item(130, 1).Value = "Tanzania"
$ws.Cells.Item(130, 2).Value = 509
$ws.Cells.Item(130, 3).Value = 0
$ws.Cells.Item(130, 4).Value = 183
$ws.Cells.Item(130, 5).Value = 305
$ws.Cells.Item(130, 6).Value = 0
$ws.Cells.Item(130, 7).Value = 0
$ws.Cells.Item(130, 8).Value = 21

$ws.Cells.Item(131, 1).Value = "Reunion"
$ws.Cells.Item(131, 2).Value = 447
$ws.Cells.Item(131, 3).Value = 1
$ws.Cells.Item(131, 4).Value = 411
$ws.Cells.Item(131, 5).Value = 35
$ws.Cells.Item(131, 6).Value = 0
$ws.Cells.Item(131, 7).Value = 1
$ws.Cells.Item(131, 8).Value = 1

$ws.Cells.Item(132, 1).Value = "Taiwan"
$ws.Cells.Item(132, 2).Value = 440
$ws.Cells.Item(132, 3).Value = 0
$ws.Cells.Item(132, 4).Value = 402
$ws.Cells.Item(132, 5).Value = 31
$ws.Cells.Item(132, 6).Value = 0
$ws.Cells.Item(132, 7).Value = 0
$ws.Cells.Item(132, 8).Value = 7

$ws.Cells.Item(133, 1).Value = "Nepal"
$ws.Cells.Item(133, 2).Value = 427
$ws.Cells.Item(133, 3).Value = 25
$ws.Cells.Item(133, 4).Value = 45
$ws.Cells.Item(133, 5).Value = 380
$ws.Cells.Item(133, 6).Value = 0
$ws.Cells.Item(133, 7).Value = 0
$ws.Cells.Item(133, 8).Value = 2

$ws.Cells.Item(134, 1).Value = "Congo"
$ws.Cells.Item(134, 2).Value = 420
$ws.Cells.Item(134, 3).Value = 0
$ws.Cells.Item(134, 4).Value = 132
$ws.Cells.Item(134, 5).Value = 273
$ws.Cells.Item(134, 6).Value = 0
$ws.Cells.Item(134, 7).Value = 0
$ws.Cells.Item(134, 8).Value = 15

$ws.Cells.Item(135, 1).Value = "Estado de Palestina"
$ws.Cells.Item(135, 2).Value = 398
$ws.Cells.Item(135, 3).Value = 7
$ws.Cells.Item(135, 4).Value = 346
$ws.Cells.Item(135, 5).Value = 50
$ws.Cells.Item(135, 6).Value = 0
$ws.Cells.Item(135, 7).Value = 0
$ws.Cells.Item(135, 8).Value = 2
